$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H64").Value = 2498.2
$ws_ALC.Range("I64").Value = 2326.25
$ws_ALC.Range("J64").Value = 2708.361
$ws_ALC.Range("K64").Value = 2326.25
$ws_ALC.Range("L64").Value = 2708.361
$ws_ALC.Range("M64").Value = -2078.25
$ws_ALC.Range("N64").Value = -3204.361
$ws_ALC.Range("H67").Value = 2498.2
$ws_ALC.Range("I67").Value = 2326.25
$ws_ALC.Range("J67").Value = 2708.361
$ws_ALC.Range("K67").Value = 2326.25
$ws_ALC.Range("L67").Value = 2708.361
$ws_ALC.Range("M67").Value = -1468.25
$ws_ALC.Range("N67").Value = -4424.361
$ws_ALC.Range("H70").Value = 2867.6667
$ws_ALC.Range("I70").Value = 2867.6667
$ws_ALC.Range("J70").Value = 0
$ws_ALC.Range("K70").Value = 8603.000100000001
$ws_ALC.Range("L70").Value = 0
$ws_ALC.Range("M70").Value = -8333.000100000001
$ws_ALC.Range("N70").ClearContents()
$ws_ALC.Range("H73").Value = 2867.6667
$ws_ALC.Range("I73").Value = 2867.6667
$ws_ALC.Range("J73").Value = 0
$ws_ALC.Range("K73").Value = 8603.000100000001
$ws_ALC.Range("L73").Value = 0
$ws_ALC.Range("M73").Value = -7667.000100000001
$ws_ALC.Range("N73").ClearContents()
$ws_ALC.Range("H86").Value = 1405883
$ws_ALC.Range("I86").Value = 43999
$ws_ALC.Range("J86").Value = 1916589.5
$ws_ALC.Range("K86").Value = 43999
$ws_ALC.Range("L86").Value = 1916589.5
$ws_ALC.Range("M86").Value = -42876
$ws_ALC.Range("N86").Value = -1918835.5
$ws_ALC.Range("H89").Value = 1405883
$ws_ALC.Range("I89").Value = 43999
$ws_ALC.Range("J89").Value = 1916589.5
$ws_ALC.Range("K89").Value = 219995
$ws_ALC.Range("L89").Value = 9582947.5
$ws_ALC.Range("M89").Value = -214379
$ws_ALC.Range("N89").Value = -9594179.5
$ws_ALC.Range("H125").Value = 1075
$ws_ALC.Range("I125").Value = 1100
$ws_ALC.Range("J125").Value = 1060
$ws_ALC.Range("K125").Value = 9900
$ws_ALC.Range("L125").Value = 9540
$ws_ALC.Range("M125").Value = -7440
$ws_ALC.Range("N125").Value = -14460
$ws_ALC.Range("H132").Value = 7817761
$ws_ALC.Range("I132").Value = 9438491
$ws_ALC.Range("J132").Value = 8788.817999999999
$ws_ALC.Range("K132").Value = 28315473
$ws_ALC.Range("L132").Value = 26366.454
$ws_ALC.Range("M132").Value = -28312943
$ws_ALC.Range("N132").Value = -31426.454
$ws_ALC.Range("H137").Value = 1193.3518
$ws_ALC.Range("I137").Value = 1035.9762
$ws_ALC.Range("J137").Value = 1744.1666
$ws_ALC.Range("K137").Value = 3107.9286
$ws_ALC.Range("L137").Value = 5232.4998
$ws_ALC.Range("M137").Value = -557.9286000000002
$ws_ALC.Range("N137").Value = -10332.4998
$ws_ALC.Range("H138").Value = 2069.681
$ws_ALC.Range("I138").Value = 1307.4482
$ws_ALC.Range("J138").Value = 3297.7222
$ws_ALC.Range("K138").Value = 3922.3446
$ws_ALC.Range("L138").Value = 9893.1666
$ws_ALC.Range("M138").Value = 1217.6554
$ws_ALC.Range("N138").Value = -20173.1666
$ws_ALC.Range("H141").Value = 2051.1042
$ws_ALC.Range("I141").Value = 1286.6586
$ws_ALC.Range("J141").Value = 6528.5713
$ws_ALC.Range("K141").Value = 3859.9758
$ws_ALC.Range("L141").Value = 19585.7139
$ws_ALC.Range("M141").Value = 1320.0242
$ws_ALC.Range("N141").Value = -29945.7139

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 1113.98
$ws_ARM.Range("I32").Value = 1081.2783
$ws_ARM.Range("J32").Value = 2171.3333
$ws_ARM.Range("K32").Value = 1081.2783
$ws_ARM.Range("L32").Value = 2171.3333
$ws_ARM.Range("M32").Value = -794.2782999999999
$ws_ARM.Range("N32").Value = -2745.3333
$ws_ARM.Range("H61").Value = 1139.5416
$ws_ARM.Range("I61").Value = 1120.6875
$ws_ARM.Range("J61").Value = 1177.25
$ws_ARM.Range("K61").Value = 1120.6875
$ws_ARM.Range("L61").Value = 1177.25
$ws_ARM.Range("M61").Value = -908.6875
$ws_ARM.Range("N61").Value = -1601.25
$ws_ARM.Range("H63").Value = 2137.8572
$ws_ARM.Range("I63").Value = 2077.5
$ws_ARM.Range("K63").Value = 2077.5
$ws_ARM.Range("M63").Value = -1391.5
$ws_ARM.Range("H66").Value = 2137.8572
$ws_ARM.Range("I66").Value = 2077.5
$ws_ARM.Range("K66").Value = 10387.5
$ws_ARM.Range("M66").Value = -6955.5
$ws_ARM.Range("H74").Value = 1546.6923
$ws_ARM.Range("I74").Value = 1738.1052
$ws_ARM.Range("J74").Value = 1027.1428
$ws_ARM.Range("K74").Value = 1738.1052
$ws_ARM.Range("L74").Value = 1027.1428
$ws_ARM.Range("M74").Value = -864.1052
$ws_ARM.Range("N74").Value = -2775.1428
$ws_ARM.Range("H77").Value = 1546.6923
$ws_ARM.Range("I77").Value = 1738.1052
$ws_ARM.Range("J77").Value = 1027.1428
$ws_ARM.Range("K77").Value = 8690.526
$ws_ARM.Range("L77").Value = 5135.714
$ws_ARM.Range("M77").Value = -4322.526
$ws_ARM.Range("N77").Value = -13871.714
$ws_ARM.Range("H136").Value = 1139.5416
$ws_ARM.Range("I136").Value = 1120.6875
$ws_ARM.Range("J136").Value = 1177.25
$ws_ARM.Range("K136").Value = 3362.0625
$ws_ARM.Range("L136").Value = 3531.75
$ws_ARM.Range("M136").Value = -812.0625
$ws_ARM.Range("N136").Value = -8631.75

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H62").Value = 4325.273
$ws_CRP.Range("I62").Value = 2866.3333
$ws_CRP.Range("J62").Value = 4872.375
$ws_CRP.Range("K62").Value = 2866.3333
$ws_CRP.Range("L62").Value = 4872.375
$ws_CRP.Range("M62").Value = -2242.3333
$ws_CRP.Range("N62").Value = -6120.375
$ws_CRP.Range("H65").Value = 4325.273
$ws_CRP.Range("I65").Value = 2866.3333
$ws_CRP.Range("J65").Value = 4872.375
$ws_CRP.Range("K65").Value = 14331.6665
$ws_CRP.Range("L65").Value = 24361.875
$ws_CRP.Range("M65").Value = -11211.6665
$ws_CRP.Range("N65").Value = -30601.875
$ws_CRP.Range("H134").Value = 1262.7188
$ws_CRP.Range("I134").Value = 1272.28
$ws_CRP.Range("J134").Value = 1228.5714
$ws_CRP.Range("K134").Value = 3816.84
$ws_CRP.Range("L134").Value = 3685.7142
$ws_CRP.Range("M134").Value = -1281.84
$ws_CRP.Range("N134").Value = -8755.7142

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H76").Value = 1500
$ws_CUL.Range("I76").Value = 1000
$ws_CUL.Range("K76").Value = 3000
$ws_CUL.Range("M76").Value = -2617
$ws_CUL.Range("H79").Value = 1500
$ws_CUL.Range("I79").Value = 1000
$ws_CUL.Range("K79").Value = 3000
$ws_CUL.Range("M79").Value = -1674
$ws_CUL.Range("H131").Value = 835.13
$ws_CUL.Range("J131").Value = 861.68475
$ws_CUL.Range("L131").Value = 2585.05425
$ws_CUL.Range("N131").Value = -12665.05425

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H132").Value = 21745484
$ws_LTW.Range("I132").Value = 41668852
$ws_LTW.Range("J132").Value = 10900.228
$ws_LTW.Range("K132").Value = 125006556
$ws_LTW.Range("L132").Value = 32700.684
$ws_LTW.Range("M132").Value = -125004026
$ws_LTW.Range("N132").Value = -37760.68399999999
